# CIMS_FIC_MB workbook updates — "industry sector updates for BC and ON -
# changes from IESO project calibration"
#
# 1. Rename the "Natural Gas Extraction" branch/sector to
#    "Natural Gas Production" everywhere it appears on the sheet
#    (Branch column A and Sector column D).
# 2. Break the stale external workbook link (to the BC calibration file)
#    and replace the broken [1]BC!$M$29 / [1]BC!$M$30 (#REF!) formulas in
#    M13:M14 with the calibrated literal value, letting the existing
#    shared formulas in N:W pick up the corrected value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Text rename across the whole sheet -------------------------------
$null = $ws.Cells.Replace("Natural Gas Extraction", "Natural Gas Production")

# --- 2. Drop the dead external reference ----------------------------------
$links = $wb.LinkSources()
foreach ($link in $links) {
  $wb.BreakLink($link, 1)
}

# --- 3. Restate the calibrated values (was #REF! via [1]BC!$M$29/$M$30) --
$ws.Range("M13").Value = 400000000
$ws.Range("M14").Value = 400000000
